# "Adding loadedViews into memory"
#
# - Adds a new "loadedView" allocation row (row 5) to the BANK61 sheet,
#   right after the existing "viewTab" row (row 4), following the same
#   layout/formula pattern as the rows above it.
# - Leaves the BANK61 sheet as the active/selected sheet with D13 selected
#   (the sheet's running-total cell), which also clears the previous
#   tab-selected state on the Golden sheet.

$wb = $excel.ActiveWorkbook

$bank61 = $wb.Worksheets.Item("BANK61")

# New row 5: loadedView segment entry
$bank61.Range("A5").Value = "loadedView"
$bank61.Range("B5").Formula = "=B4+ E4+1"
$bank61.Range("C5").Value = 6
$bank61.Range("D5").Value = 256
$bank61.Range("E5").Formula = "=C5*D5"

# Make BANK61 the active sheet/tab and select D13
$bank61.Activate() | Out-Null
$bank61.Range("D13").Select() | Out-Null
